# GSYE-681: Added support for SCMStorage strategy in the CDS.
# Modified schema in order to support the SCMStorage strategy.
#
# The "Storage" sheet's per-battery spec columns (Capacity [kWh],
# Minimum allowed SoC [-], Maximum power [kW]) are no longer fixed
# columns in the schema -- clear the header labels and the sample
# values for the two existing batteries, keeping their styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Storage")

# Clear the header (row 1) and the two battery rows (rows 2-3) for
# columns C:E -- this drops "Capacity [kWh]", "Minimum allowed SoC [-]"
# and "Maximum power [kW]" plus their sample data, while leaving the
# cell formatting (style) in place.
$null = $ws.Range("C1:E3").ClearContents()

# The Storage sheet becomes the active tab/selection of the workbook.
$null = $ws.Activate()
$null = $ws.Range("C1:E8").Select()
